# BIS-1002: removed "Internal Assignment" column from export.
#
# The "Internal Assignment" property column (column O) is no longer part of
# the exported sample-type sheet. Its header cell (O4) and the per-property
# "FALSE" values below it (O5:O9) are cleared out, while keeping their
# existing cell styles/formatting intact. Clearing (rather than deleting the
# column) also drops the now-unreferenced "Internal Assignment" shared
# string from the workbook.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear the "Internal Assignment" header cell and its column values, keeping
# cell formatting (style) untouched.
$ws.Range("O4").ClearContents()
$ws.Range("O5:O9").ClearContents()

# Reflect the edit in the current selection, same as Excel would leave it
# after clearing that range interactively.
$ws.Range("O4:O9").Select()
